$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    "movie_name",
    "stimuli/duck_video.mp4",
    "stimuli/Jester.mp4",
    "stimuli/Kids_Jack_01.mp4",
    "stimuli/Kids_Jack_02.mp4",
    "stimuli/Kids_Jack_03.mp4",
    "stimuli/Kids_Jack_04.mp4",
    "stimuli/Kids_Jack_05.mp4",
    "stimuli/Kids_Jack_06.mp4",
    "stimuli/Kids_Jack_07.mp4",
    "stimuli/Kids_Jack_08.mp4",
    "stimuli/Kids_Jack_09.mp4",
    "stimuli/Kids_Jack_10.mp4",
    "stimuli/Kids_Jack_11.mp4",
    "stimuli/Kids_Jack_12.mp4",
    "stimuli/Kids_Jack_13.mp4",
    "stimuli/Kids_Jack_14.mp4",
    "stimuli/Kids_Jack_15.mp4",
    "stimuli/Kids_Jack_16.mp4",
    "stimuli/Kids_Jack_17.mp4",
    "stimuli/Kids_Jack_18.mp4",
    "stimuli/Kids_Jack_19.mp4",
    "stimuli/Kids_Jack_20.mp4",
    "stimuli/Kids_Jack_21.mp4",
    "stimuli/Kids_Jack_22.mp4",
    "stimuli/Kids_Jack_23.mp4",
    "stimuli/Kids_Jack_24.mp4",
    "stimuli/Kids_Jack_25.mp4",
    "stimuli/Kids_Jack_26.mp4",
    "stimuli/Kids_Jack_27.mp4",
    "stimuli/Kids_Jack_28.mp4",
    "stimuli/Kids_Jack_29.mp4",
    "stimuli/P01.mp4",
    "stimuli/P02s.mp4",
    "stimuli/P03.mp4",
    "stimuli/P04.mp4",
    "stimuli/P05.mp4",
    "stimuli/P06.mp4",
    "stimuli/P07s.mp4",
    "stimuli/P08.mp4",
    "stimuli/P08.mpg",
    "stimuli/P09.mp4",
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("D12").Select()
